# Commit 07032018 at 3.20PM
# The "InputData" sheet's Browser test-data cell (B2) is changed from "IE"
# to "Firefox", and the sheet's active selection moves from D7 to B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 held the Browser value "IE" -> update it to "Firefox"
$ws.Range("B2").Value = "Firefox"

# Active selection moves to B3
$ws.Range("B3").Select()
